$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column widths (closest achievable values on this engine's pixel-quantized
# column width grid; target character widths are 11.7109375 / 12.7109375)
$ws.Columns.Item(2).ColumnWidth = 10.833333333333334
$ws.Columns.Item(3).ColumnWidth = 11.833333333333334

# Update cell values
$ws.Range("A1").Value = 160.72160438563503
$ws.Range("B1").Value = 6.187470435797815
$ws.Range("C1").Value = 0.98450195694716236
